$d = $word.ActiveDocument

$newText = "年キャンペーン期間 対象：Bootes: 5月14〜23日、6月13〜22日、7月12〜21日"

# Locate every paragraph whose text still contains the old, multi-run
# "2018..." campaign-period sentence (there are 4 occurrences in this
# document) and replace its whole content with a single plain run.
$targets = @()
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t.StartsWith("201")) {
        $targets += $i
    }
}

# Walk backwards so earlier paragraph indices stay valid while we edit.
for ($k = $targets.Count - 1; $k -ge 0; $k--) {
    $idx = $targets[$k]
    $p = $d.Paragraphs($idx)
    $start = $p.Range.Start
    $end = $p.Range.End - 1   # exclude the paragraph mark

    $r = $d.Range($start, $end)
    $r.Delete()

    $r2 = $d.Range($start, $start)
    $r2.InsertAfter($newText)
}
